# Update "报名/观展人数" (column F) figures across sheets, per gh-pages
# regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4940
$ws1.Range("F8").Value  = 89
$ws1.Range("F12").Value = 361
$ws1.Range("F19").Value = 70
$ws1.Range("F25").Value = 544
$ws1.Range("F26").Value = 1022
$ws1.Range("F28").Value = 1917
$ws1.Range("F29").Value = 2376
$ws1.Range("F30").Value = 1156
$ws1.Range("F32").Value = 84
$ws1.Range("F33").Value = 320
$ws1.Range("F34").Value = 330
$ws1.Range("F36").Value = 657
$ws1.Range("F39").Value = 712
$ws1.Range("F42").Value = 580
$ws1.Range("F43").Value = 264

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 212

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 4940
$ws4.Range("F15").Value = 361
$ws4.Range("F26").Value = 1022
$ws4.Range("F28").Value = 1917
$ws4.Range("F29").Value = 2376
$ws4.Range("F31").Value = 1156
$ws4.Range("F35").Value = 84
$ws4.Range("F36").Value = 320
$ws4.Range("F37").Value = 331
$ws4.Range("F41").Value = 657
$ws4.Range("F42").Value = 712
$ws4.Range("F45").Value = 580
$ws4.Range("F46").Value = 264
